$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2023-09-25 Monday" "2023-09-26 Tuesday"

Replace-Text "93×54=5022" "11×57=627"
Replace-Text "94×43=4042" "84×43=3612"
Replace-Text "44×58=2552" "71×15=1065"
Replace-Text "12×55=660" "96×62=5952"
Replace-Text "55×36=1980" "49×17=833"
Replace-Text "40×38=1520" "55×36=1980"
Replace-Text "18×43=774" "58×11=638"
Replace-Text "56×36=2016" "48×56=2688"
Replace-Text "77×96=7392" "71×74=5254"
Replace-Text "40×27=1080" "54×49=2646"

Replace-Text "95×32=3040" "41×16=656"
Replace-Text "89×44=3916" "71×34=2414"
Replace-Text "49×46=2254" "42×18=756"
Replace-Text "26×37=962" "44×50=2200"
Replace-Text "69×65=4485" "54×39=2106"

Replace-Text "84×44=3696" "87×92=8004"
Replace-Text "51×64=3264" "31×80=2480"
Replace-Text "59×98=5782" "51×65=3315"
Replace-Text "92×57=5244" "84×23=1932"
Replace-Text "20×70=1400" "24×29=696"

Replace-Text "65×93=6045" "66×67=4422"
Replace-Text "91×66=6006" "31×35=1085"
Replace-Text "50×28=1400" "42×65=2730"
Replace-Text "16×30=480" "26×94=2444"
Replace-Text "58×86=4988" "68×82=5576"
